# Daily attendance processing - 2025-10-21 07:20:24
# Normalizes the "Recorded By" column (G) so that entries logged as
# "System, <email>[, <email2>]" have the leading "System" marker moved
# to the end of the list (capitalized), matching the latest export format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -notlike "System, *") { continue }
    if ($val -like "*admin@admin.com*") { continue }

    $parts = $val -split ", "
    if ($parts.Count -lt 2) { continue }

    $last = $parts[$parts.Count - 1]

    $parts[0] = $last
    $parts[$parts.Count - 1] = "System"

    $cell.Value2 = [string]::Join(", ", $parts)
}
